$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the "You parents" -> "Your parents" typo in the second paragraph,
#    and make sure the cursor's last-edit position (the "_GoBack" bookmark)
#    ends up right after the inserted "r", exactly as Word would leave it
#    after a live edit.
# ---------------------------------------------------------------------------

# Locate the exact text that needs fixing instead of hard-coding offsets.
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Execute("You parents have said you can have some friends over ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$runStart = $find.Start
$runEnd = $find.End

# "You" occupies the first three characters of that run.
$youEnd = $runStart + 3

# Protective, zero-width bookmark right at the end of the run being edited.
# It has no semantic purpose of its own (it's removed below); its only job
# is to stop the engine's run-reflow from bleeding into the unrelated runs
# that follow (e.g. "during the school holidays ...") when we edit the text
# just before it.
$d.Bookmarks.Add("ZZtmpGuard", $d.Range($runEnd, $runEnd)) | Out-Null

# Insert the missing "r" right after "You" to fix the typo.
$d.Range($youEnd, $youEnd).Text = "r"
$rEnd = $youEnd + 1

# Split "You" into its own run (no bookmark left behind - add then delete).
$d.Bookmarks.Add("ZZtmpSplit", $d.Range($runStart, $youEnd)) | Out-Null
$d.Bookmarks("ZZtmpSplit").Delete()

# Drop a collapsed "_GoBack" bookmark immediately after the "r" - this both
# marks the last-edit position and forces the run split between "r" and
# " parents have said ...". Word only ever keeps a single "_GoBack" bookmark,
# so this naturally supersedes/moves the one that used to sit elsewhere.
$d.Bookmarks.Add("_GoBack", $d.Range($rEnd, $rEnd)) | Out-Null

# Remove the temporary guard bookmark; it has done its job.
$d.Bookmarks("ZZtmpGuard").Delete()

# ---------------------------------------------------------------------------
# 2) Make sure no stray "_GoBack" bookmark remains after "Pseudocode" (it
#    should already have been superseded above, but double-check/clean up
#    defensively in case of a different engine behaviour).
# ---------------------------------------------------------------------------

$pseudoRange = $d.Content.Find
$pseudoHit = $d.Content
$pseudoHit.Find.Execute("Pseudocode", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

if ($pseudoHit.Find.Found) {
    $after = $d.Range($pseudoHit.End, $pseudoHit.End)
    if ($d.Bookmarks.Exists("_GoBack")) {
        $gb = $d.Bookmarks("_GoBack")
        if ($gb.Start -ge $pseudoHit.Start -and $gb.Start -le ($pseudoHit.End + 2)) {
            $gb.Delete()
        }
    }
}
